$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the text header labels in row 1 (B1:L1) with plain numeric
# index values (0..10), keeping existing cell formatting (bold/center/border).
$ws.Range("B1").Value = 0
$ws.Range("C1").Value = 1
$ws.Range("D1").Value = 2
$ws.Range("E1").Value = 3
$ws.Range("F1").Value = 4
$ws.Range("G1").Value = 5
$ws.Range("H1").Value = 6
$ws.Range("I1").Value = 7
$ws.Range("J1").Value = 8
$ws.Range("K1").Value = 9
$ws.Range("L1").Value = 10
